$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.0999
$ws.Range("E2").Value = 0.103
$ws.Range("F2").Value = 0.219
$ws.Range("G2").Value = 0.09150696864111497
$ws.Range("H2").Value = 0.09150696864111497
$ws.Range("I2").Value = 0.08997055783743976
$ws.Range("J2").Value = 0.07371135699102971
$ws.Range("K2").Value = 92.41
$ws.Range("L2").Value = 0.04024825783972125
$ws.Range("M2").Value = 68.33
$ws.Range("N2").Value = 0.0269812438302073
$ws.Range("O2").Value = 0.7394221404609891
$ws.Range("P2").Value = 58.23
$ws.Range("Q2").Value = 0.02299308983218164
$ws.Range("R2").Value = 0.6301266096742777
$ws.Range("S2").Value = 10.1
$ws.Range("T2").Value = 0.1478120883945558
$ws.Range("U2").Value = 190.4
$ws.Range("V2").Value = 0.07518262586377097
$ws.Range("W2").Value = 0.09054054054054055
$ws.Range("X2").Value = 0.06353781670736119
$ws.Range("Y2").Value = 0.02700272383317935
$ws.Range("Z2").Value = 2.381772397026423
$ws.Range("AA2").Value = 0.1334230493915533
$ws.Range("AB2").Value = 0.06353781670736119
$ws.Range("AC2").Value = 0.06988523268419214
$ws.Range("AD2").Value = 30.4
$ws.Range("AE2").Value = 11.78799602619153
$ws.Range("AF2").Value = 42.18799602619153
$ws.Range("AG2").Value = -148.2120039738085
$ws.Range("AH2").Value = 0.01638567317333403
$ws.Range("AI2").Value = 0.03142523147400154
$ws.Range("AJ2").Value = -0.06216195535976701
$ws.Range("AK2").Value = -0.1286464267356528
$ws.Range("AL2").Value = 47.125
$ws.Range("AM2").Value = 47.125
$ws.Range("AN2").Value = 0.1401373715023279
$ws.Range("AO2").Value = 4.350132625994695
$ws.Range("AP2").Value = -0.683225021775727
$ws.Range("AQ2").Value = 4.350132625994695

# Row 3
$ws.Range("D3").Value = 0.186
$ws.Range("E3").Value = 0.00446
$ws.Range("G3").Value = 0.081436680991034
$ws.Range("H3").Value = 0.081436680991034
$ws.Range("I3").Value = 0.08062624054048581
$ws.Range("J3").Value = 0.0670520113638583
$ws.Range("K3").Value = 55.1
$ws.Range("L3").Value = 0.02923232001697703
$ws.Range("M3").Value = 31.2
$ws.Range("N3").Value = 0.01480989224854037
$ws.Range("O3").Value = 0.5662431941923775
$ws.Range("P3").Value = 31.2
$ws.Range("Q3").Value = 0.01480989224854037
$ws.Range("R3").Value = 0.5662431941923775
$ws.Range("U3").Value = 137
$ws.Range("V3").Value = 0.06503061660416766
$ws.Range("W3").Value = 0.08178714561377468
$ws.Range("X3").Value = 0.06440673389375487
$ws.Range("Y3").Value = 0.01738041172001981
$ws.Range("Z3").Value = 2.987220061032287
$ws.Range("AA3").Value = 0.2002991134786823
$ws.Range("AB3").Value = 0.06394343828022246
$ws.Range("AC3").Value = 0.1363556751984599
$ws.Range("AD3").Value = 30.4
$ws.Range("AE3").Value = 11.78799602619153
$ws.Range("AF3").Value = 42.18799602619153
$ws.Range("AG3").Value = -94.81200397380847
$ws.Range("AH3").Value = 0.01963247786958056
$ws.Range("AI3").Value = 0.044667583075372
$ws.Range("AJ3").Value = -0.047125885815253
$ws.Range("AK3").Value = -0.1174159918666248
$ws.Range("AL3").Value = 46.6
$ws.Range("AM3").Value = 46.6
$ws.Range("AN3").Value = 0.1907989706897634
$ws.Range("AO3").Value = 3.227467811158798
$ws.Range("AP3").Value = -0.5950668673433029
$ws.Range("AQ3").Value = 3.227467811158798

# Row 4
$ws.Range("D4").Value = 0.0333
$ws.Range("E4").Value = 0.143
$ws.Range("G4").Value = 0.1835548172757475
$ws.Range("H4").Value = 0.1835548172757475
$ws.Range("I4").Value = 0.1719269102990033
$ws.Range("J4").Value = 0.1371354016794412
$ws.Range("K4").Value = 28.6
$ws.Range("L4").Value = 0.1187707641196013
$ws.Range("M4").Value = 32.4
$ws.Range("N4").Value = 0.106020942408377
$ws.Range("O4").Value = 1.132867132867133
$ws.Range("P4").Value = 22.3
$ws.Range("Q4").Value = 0.07297120418848167
$ws.Range("R4").Value = 0.7797202797202797
$ws.Range("S4").Value = 10.1
$ws.Range("T4").Value = 0.3117283950617283
$ws.Range("U4").Value = 40.5
$ws.Range("V4").Value = 0.1325261780104712
$ws.Range("W4").Value = 0.09688346883468836
$ws.Range("X4").Value = 0.06353781670736119
$ws.Range("Y4").Value = 0.03334565212732717
$ws.Range("Z4").Value = 0.972929292929293
$ws.Range("AA4").Value = 0.1334230493915533
$ws.Range("AB4").Value = 0.06353781670736119
$ws.Range("AC4").Value = 0.06988523268419214
$ws.Range("AG4").Value = -40.5
$ws.Range("AJ4").Value = -0.1527725386646548
$ws.Range("AK4").Value = -0.1553509781357883
$ws.Range("AL4").Value = 0.041
$ws.Range("AM4").Value = 0.041
$ws.Range("AO4").Value = 1009.756097560975
$ws.Range("AP4").Value = -0.9267734553775743
$ws.Range("AQ4").Value = 1009.756097560975

# Row 5
$ws.Range("D5").Value = 0.0999
$ws.Range("E5").Value = 0.103
$ws.Range("F5").Value = 0.219
$ws.Range("G5").Value = 0.0728126834997064
$ws.Range("H5").Value = 0.0728126834997064
$ws.Range("I5").Value = 0.07751027598355842
$ws.Range("J5").Value = 0.06422280010066268
$ws.Range("K5").Value = 8.710000000000001
$ws.Range("L5").Value = 0.05114503816793893
$ws.Range("M5").Value = 4.73
$ws.Range("N5").Value = 0.03935108153078203
$ws.Range("O5").Value = 0.5430539609644087
$ws.Range("P5").Value = 4.73
$ws.Range("Q5").Value = 0.03935108153078203
$ws.Range("R5").Value = 0.5430539609644087
$ws.Range("U5").Value = 12.9
$ws.Range("V5").Value = 0.1073211314475874
$ws.Range("W5").Value = 0.09054054054054055
$ws.Range("X5").Value = 0.06353781670736119
$ws.Range("Y5").Value = 0.02700272383317935
$ws.Range("Z5").Value = 1.991812865497076
$ws.Range("AA5").Value = 0.1279197994987469
$ws.Range("AB5").Value = 0.06353781670736119
$ws.Range("AC5").Value = 0.06438198279138566
$ws.Range("AG5").Value = -12.9
$ws.Range("AJ5").Value = -0.1202236719478099
$ws.Range("AK5").Value = -0.1537544696066746
$ws.Range("AL5").Value = 0.484
$ws.Range("AM5").Value = 0.484
$ws.Range("AO5").Value = 27.27272727272727
$ws.Range("AP5").Value = -0.9280575539568345
$ws.Range("AQ5").Value = 27.27272727272727